$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 156-187 (cell-level changes) ---
# Row 156
$ws.Range("D156").Value = 44995
$ws.Range("K156").Value = 800
$ws.Range("M156").Value = 800
$ws.Range("P156").Value = 800

# Row 157
$ws.Range("D157").Value = 44995

# Row 158
$ws.Range("D158").Value = 44798
$ws.Range("K158").Value = 700
$ws.Range("L158").Value = 800

# Row 159
$ws.Range("D159").Value = 44798
$ws.Range("K159").Value = 600
$ws.Range("L159").Value = 600
$ws.Range("M159").Value = 600
$ws.Range("P159").Value = 600

# Row 160
$ws.Range("D160").Value = 44974
$ws.Range("J160").Value = 200
$ws.Range("K160").Value = 750
$ws.Range("L160").Value = 750
$ws.Range("M160").Value = 750
$ws.Range("P160").Value = 750

# Row 161
$ws.Range("D161").Value = 44974
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 150
$ws.Range("K161").Value = 650
$ws.Range("L161").Value = 650
$ws.Range("M161").Value = 650
$ws.Range("P161").Value = 650

# Row 162
$ws.Range("D162").Value = 44963
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 300
$ws.Range("K162").Value = 800
$ws.Range("L162").Value = 900
$ws.Range("M162").Value = 850
$ws.Range("P162").Value = 850

# Row 163
$ws.Range("D163").Value = 44966
$ws.Range("J163").Value = 300
$ws.Range("K163").Value = 800
$ws.Range("M163").Value = 800
$ws.Range("P163").Value = 800

# Row 164
$ws.Range("D164").Value = 44966
$ws.Range("J164").Value = 200

# Row 165
$ws.Range("D165").Value = 44763
$ws.Range("K165").Value = 700
$ws.Range("L165").Value = 800
$ws.Range("M165").Value = 750
$ws.Range("P165").Value = 750

# Row 166
$ws.Range("D166").Value = 44763
$ws.Range("I166").Value = "Segunda"
$ws.Range("J166").Value = 150
$ws.Range("K166").Value = 600
$ws.Range("M166").Value = 600
$ws.Range("P166").Value = 600

# Row 167
$ws.Range("D167").Value = 44701
$ws.Range("J167").Value = 200
$ws.Range("K167").Value = 550
$ws.Range("L167").Value = 600
$ws.Range("M167").Value = 575
$ws.Range("O167").Value = "Provincia de Diguillín"
$ws.Range("P167").Value = 575

# Row 168
$ws.Range("D168").Value = 44608
$ws.Range("J168").Value = 160
$ws.Range("K168").Value = 550
$ws.Range("L168").Value = 600
$ws.Range("M168").Value = 575
$ws.Range("P168").Value = 575

# Row 169
$ws.Range("D169").Value = 44208
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 130
$ws.Range("K169").Value = 1800
$ws.Range("L169").Value = 2000
$ws.Range("M169").Value = 1908
$ws.Range("O169").Value = "Provincia de Cautín"
$ws.Range("P169").Value = 1908

# Row 170
$ws.Range("D170").Value = 44839
$ws.Range("J170").Value = 200
$ws.Range("K170").Value = 700
$ws.Range("L170").Value = 800
$ws.Range("M170").Value = 750
$ws.Range("P170").Value = 750

# Row 171
$ws.Range("D171").Value = 44839
$ws.Range("I171").Value = "Segunda"
$ws.Range("J171").Value = 250
$ws.Range("K171").Value = 600
$ws.Range("L171").Value = 600
$ws.Range("M171").Value = 600
$ws.Range("P171").Value = 600

# Row 172
$ws.Range("D172").Value = 44637
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 120
$ws.Range("K172").Value = 550
$ws.Range("M172").Value = 575
$ws.Range("P172").Value = 575

# Row 173
$ws.Range("D173").Value = 44771
$ws.Range("J173").Value = 200
$ws.Range("K173").Value = 700
$ws.Range("L173").Value = 800
$ws.Range("M173").Value = 750
$ws.Range("N173").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("P173").Value = 750
$ws.Range("Q173").Value = 1

# Row 174
$ws.Range("D174").Value = 44771
$ws.Range("I174").Value = "Segunda"
$ws.Range("K174").Value = 600
$ws.Range("M174").Value = 600
$ws.Range("P174").Value = 600

# Row 175
$ws.Range("D175").Value = 44160
$ws.Range("J175").Value = 190
$ws.Range("K175").Value = 1300
$ws.Range("L175").Value = 1500
$ws.Range("M175").Value = 1395
$ws.Range("N175").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("P175").Value = 930
$ws.Range("Q175").Value = 1.5

# Row 176
$ws.Range("D176").Value = 44665
$ws.Range("J176").Value = 200
$ws.Range("K176").Value = 550
$ws.Range("L176").Value = 600
$ws.Range("M176").Value = 575
$ws.Range("P176").Value = 575

# Row 177
$ws.Range("D177").Value = 44270
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 260
$ws.Range("K177").Value = 1800
$ws.Range("L177").Value = 2000
$ws.Range("M177").Value = 1908
$ws.Range("P177").Value = 1908

# Row 178
$ws.Range("D178").Value = 44882
$ws.Range("J178").Value = 400
$ws.Range("K178").Value = 600
$ws.Range("L178").Value = 700
$ws.Range("M178").Value = 650
$ws.Range("P178").Value = 650

# Row 179
$ws.Range("D179").Value = 44882
$ws.Range("I179").Value = "Segunda"
$ws.Range("J179").Value = 300
$ws.Range("K179").Value = 500
$ws.Range("L179").Value = 500
$ws.Range("M179").Value = 500
$ws.Range("P179").Value = 500

# Row 180
$ws.Range("D180").Value = 44273
$ws.Range("J180").Value = 140
$ws.Range("K180").Value = 1800
$ws.Range("L180").Value = 2000
$ws.Range("M180").Value = 1914
$ws.Range("P180").Value = 1914

# Row 181
$ws.Range("D181").Value = 44670
$ws.Range("J181").Value = 120
$ws.Range("K181").Value = 550
$ws.Range("L181").Value = 600
$ws.Range("M181").Value = 575
$ws.Range("P181").Value = 575

# Row 182
$ws.Range("D182").Value = 44610
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 100
$ws.Range("K182").Value = 550
$ws.Range("L182").Value = 600
$ws.Range("M182").Value = 575
$ws.Range("P182").Value = 575

# Row 183
$ws.Range("D183").Value = 44819
$ws.Range("J183").Value = 200
$ws.Range("K183").Value = 800
$ws.Range("L183").Value = 900
$ws.Range("M183").Value = 850
$ws.Range("P183").Value = 850

# Row 184
$ws.Range("D184").Value = 44819
$ws.Range("J184").Value = 150
$ws.Range("K184").Value = 700
$ws.Range("L184").Value = 700
$ws.Range("M184").Value = 700
$ws.Range("P184").Value = 700

# Row 185
$ws.Range("D185").Value = 44876

# Row 186
$ws.Range("D186").Value = 44876

# Row 187
$ws.Range("D187").Value = 44859
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 400
$ws.Range("L187").Value = 700
$ws.Range("M187").Value = 650
$ws.Range("P187").Value = 650

# --- Add new rows 188-189 (full new rows) ---
# Row 188
$ws.Range("A188").Value = 7
$ws.Range("B188").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C188").Value = "Ñuble"
$ws.Range("D188").Value = 44859
$ws.Range("E188").Value = 16
$ws.Range("F188").Value = 100112040
$ws.Range("G188").Value = "Cilantro"
$ws.Range("H188").Value = "Sin especificar"
$ws.Range("I188").Value = "Segunda"
$ws.Range("J188").Value = 300
$ws.Range("K188").Value = 500
$ws.Range("L188").Value = 500
$ws.Range("M188").Value = 500
$ws.Range("N188").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O188").Value = "Provincia de Diguillín"
$ws.Range("P188").Value = 500
$ws.Range("Q188").Value = 1
$ws.Range("R188").Value = "Hortaliza"
$ws.Range("D188").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 189
$ws.Range("A189").Value = 7
$ws.Range("B189").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C189").Value = "Ñuble"
$ws.Range("D189").Value = 44799
$ws.Range("E189").Value = 16
$ws.Range("F189").Value = 100112040
$ws.Range("G189").Value = "Cilantro"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Segunda"
$ws.Range("J189").Value = 100
$ws.Range("K189").Value = 600
$ws.Range("L189").Value = 600
$ws.Range("M189").Value = 600
$ws.Range("N189").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O189").Value = "Provincia de Diguillín"
$ws.Range("P189").Value = 600
$ws.Range("Q189").Value = 1
$ws.Range("R189").Value = "Hortaliza"
$ws.Range("D189").NumberFormat = "YYYY-MM-DD HH:MM:SS"
